$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve column D (Price) as text so values like "1.012" are not
# auto-converted to numbers by Excel's COM value parser.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.993.81"
$ws.Range("E2").Value = "  -1.30%  "

$ws.Range("D3").Value = "1.819.69"
$ws.Range("E3").Value = "  -1.00%  "

$ws.Range("D4").Value = "1.012"
$ws.Range("E4").Value = "  -0.29%  "

$ws.Range("D5").Value = "310.87"
$ws.Range("E5").Value = "  -1.27%  "

$ws.Range("E6").Value = "  -0.27%  "

$ws.Range("D7").Value = "0.4628"
$ws.Range("E7").Value = "  -2.30%  "

$ws.Range("D8").Value = "0.3623"
$ws.Range("E8").Value = "  -2.06%  "

$ws.Range("D9").Value = "0.07295"
$ws.Range("E9").Value = "  -2.23%  "

$ws.Range("D10").Value = "0.8648"
$ws.Range("E10").Value = "  -2.30%  "

$ws.Range("D11").Value = "19.72"
$ws.Range("E11").Value = "  -3.76%  "

$ws.Range("D12").Value = "1.873.56"
$ws.Range("E12").Value = "  +1.69%  "

$ws.Range("D13").Value = "0.07598"
$ws.Range("E13").Value = "  +3.12%  "

$ws.Range("D14").Value = "92.87"
$ws.Range("E14").Value = "  -0.39%  "

$ws.Range("D15").Value = "5.316"
$ws.Range("E15").Value = "  -3.02%  "

$ws.Range("D16").Value = "6.439"
$ws.Range("E16").Value = "  -2.12%  "

$ws.Range("D17").Value = "1.011"
$ws.Range("E17").Value = "  -0.36%  "

$ws.Range("D18").Value = "0.000008621"
$ws.Range("E18").Value = "  -2.43%  "

$ws.Range("E19").Value = "  -0.26%  "

$ws.Range("D20").Value = "27.247.10"
$ws.Range("E20").Value = "  -0.46%  "

$ws.Range("D21").Value = "14.42"
$ws.Range("E21").Value = "  -2.81%  "

$ws.Range("D22").Value = "5.157"
$ws.Range("E22").Value = "  -3.72%  "

$ws.Range("D23").Value = "10.55"
$ws.Range("E23").Value = "  -1.40%  "

$ws.Range("D24").Value = "2.097.93"
$ws.Range("E24").Value = "  +1.38%  "

$ws.Range("D25").Value = "151.22"
$ws.Range("E25").Value = "  -0.65%  "

$ws.Range("E26").Value = "  -2.60%  "

$ws.Range("D27").Value = "18.21"
$ws.Range("E27").Value = "  -2.35%  "

$ws.Range("D28").Value = "2.095"
$ws.Range("E28").Value = "  -3.25%  "

$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "115.51"
$ws.Range("E29").Value = "  -2.05%  "

$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "5.065"
$ws.Range("E30").Value = "  -3.67%  "

$ws.Range("D31").Value = "0.08882"
$ws.Range("E31").Value = "  -0.94%  "

$ws.Range("D32").Value = "2.957"
$ws.Range("E32").Value = "  +0.45%  "

$ws.Range("D33").Value = "0.7264"
$ws.Range("E33").Value = "  -4.40%  "

$ws.Range("D34").Value = "1.134"
$ws.Range("E34").Value = "  -3.77%  "

$ws.Range("D35").Value = "4.404"
$ws.Range("E35").Value = "  -3.40%  "

$ws.Range("D36").Value = "1.011"
$ws.Range("E36").Value = "  -0.26%  "

$ws.Range("D37").Value = "2.524"
$ws.Range("E37").Value = "  +6.49%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "1.072"
$ws.Range("E38").Value = "  -3.16%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "0.05243"
$ws.Range("E39").Value = "  -2.49%  "

$ws.Range("D40").Value = "0.01914"
$ws.Range("E40").Value = "  -2.47%  "

$ws.Range("D41").Value = "2.922"
$ws.Range("E41").Value = "  -2.82%  "

$ws.Range("D42").Value = "7.113"
$ws.Range("E42").Value = "  -2.29%  "

$ws.Range("D43").Value = "0.5192"
$ws.Range("E43").Value = "  -2.96%  "

$ws.Range("D44").Value = "0.1626"
$ws.Range("E44").Value = "  -2.24%  "

$ws.Range("D45").Value = "8.205"
$ws.Range("E45").Value = "  -3.97%  "

$ws.Range("D46").Value = "0.4845"
$ws.Range("E46").Value = "  -2.53%  "

$ws.Range("D47").Value = "1.011"
$ws.Range("E47").Value = "  -0.38%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "10.11"
$ws.Range("E48").Value = "  -3.72%  "

$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "103.04"
$ws.Range("E49").Value = "  -2.00%  "

$ws.Range("D50").Value = "1.633"
$ws.Range("E50").Value = "  -2.85%  "

$ws.Range("D51").Value = "0.06237"
$ws.Range("E51").Value = "  -1.36%  "
